# "Validator in Create-K. gesetzt"
#
# Applies the to-do list updates for the "Tabelle1" (sheet1), "03.10."
# (sheet2) and "Notizen Präsi" (sheet3) worksheets:
#  - several open tasks were checked off / reworded in the Tabelle1 todo
#    column (A4:A23) and the "Naechste Schritte" column (H2:H6)
#  - a few finished tasks were promoted into the "done" list (A41:A43)
#  - the little create-form field labels (E14:E18, E22:E23) got the
#    highlight style that the rest of that block already had
#  - "Notizen Praesi" (sheet3) got one more bullet appended (A18)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws3 = $wb.Worksheets.Item("Notizen Präsi")

# ---------------------------------------------------------------------------
# Tabelle1 - "Naechste Schritte" box (column H)
# ---------------------------------------------------------------------------
$ws1.Range("F2").ClearContents()

$ws1.Range("H3").Value = "create invalid!"
$ws1.Range("H3").Copy()
$ws1.Range("H3").PasteSpecial(-4122) | Out-Null

$ws1.Range("H4").Value = "done button"
$ws1.Range("H5").Value = "create Toast!"
$ws1.Range("H6").Value = "delete Modal!"

$ws1.Range("H9").Copy()
$ws1.Range("H10").PasteSpecial(-4122) | Out-Null
$ws1.Range("H10").Value = "Deployment"

# ---------------------------------------------------------------------------
# Tabelle1 - open to-do list (column A, rows 4-23)
# ---------------------------------------------------------------------------
$ws1.Range("A6").Value = "create Toast!"
$ws1.Range("A9").Value = "create new erstellen!"
$ws1.Range("A10").Value = "Suchleiste?"
$ws1.Range("A11").Value = "Deployment (Frontend über github, Backend vercel oder render)"
$ws1.Range("A12").ClearContents()
$ws1.Range("A13").Value = "Modal delete funktioniert nicht"
$ws1.Range("A14").ClearContents()
$ws1.Range("A15").Value = "code englisch"
$ws1.Range("A16").Value = "seite deutsch"
$ws1.Range("A17").Value = "Datepicker nach oben / unten begrenzen"
$ws1.Range("A18").Value = "validator bei update"
$ws1.Range("A19").ClearContents()
$ws1.Range("A20").ClearContents()
$ws1.Range("A21").ClearContents()
$ws1.Range("A22").ClearContents()
$ws1.Range("A23").ClearContents()

# form-field labels that now share the highlighted style used elsewhere
$ws1.Range("E9").Copy()
$ws1.Range("E14").PasteSpecial(-4122) | Out-Null
$ws1.Range("E15").PasteSpecial(-4122) | Out-Null
$ws1.Range("E16").PasteSpecial(-4122) | Out-Null
$ws1.Range("E17").PasteSpecial(-4122) | Out-Null
$ws1.Range("E18").PasteSpecial(-4122) | Out-Null
$ws1.Range("E22").PasteSpecial(-4122) | Out-Null
$ws1.Range("E23").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# Tabelle1 - "erledigt" / done list (column A, rows 28-40) grows by three
# ---------------------------------------------------------------------------
$ws1.Range("A40").Copy()
$ws1.Range("A41").PasteSpecial(-4122) | Out-Null
$ws1.Range("A42").PasteSpecial(-4122) | Out-Null
$ws1.Range("A43").PasteSpecial(-4122) | Out-Null

$ws1.Range("A41").Value = "Backend und Frontend verknüpfen"
$ws1.Range("A42").Value = "update Methode hinterlegen"
$ws1.Range("A43").Value = "create invalid/valid !!"

$ws1.Range("E14").Select()

# ---------------------------------------------------------------------------
# Notizen Präsi - one more note appended
# ---------------------------------------------------------------------------
$ws3.Range("A18").Value = "unterschied patch und put"
$ws3.Application.ActiveWindow.ScrollRow = 10
$ws3.Range("A18").Select()
